# Daily auto push: insert a new data row for 2026/01/18 06:46 UTC run.
# A new row is inserted before the existing row 648, shifting all
# subsequent rows down by one (old row 648 becomes row 649, etc.),
# and the worksheet's used-range dimension grows from D689 to D690.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 648 (existing rows 648.. shift down).
$ws.Rows.Item(648).Insert()

# Force column A to be stored as text so the date string is not
# reinterpreted as a numeric date serial value.
$ws.Range("A648").NumberFormat = "@"
$ws.Range("A648").Value = "2026/01/18"
$ws.Range("B648").Value = "日"
$ws.Range("C648").Value = 13
$ws.Range("D648").Value = 23

# Drop any formatting picked up from the NumberFormat change / insert
# so the new row matches the plain (unstyled) look of the other data rows.
$ws.Range("A648:D648").ClearFormats()
